# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" (shared string used
#    by Overview!E2/F2, zh-cn!C2, de-de!C2 — update every occurrence so the
#    old string becomes unused and is dropped, and the new one is written).
# 2) Narrow the "Status" column: Overview columns E & F, and column C on the
#    zh-cn / de-de sheets (17.22 chars -> 13.41 chars, i.e. ColumnWidth 12.5
#    in Excel's character-width units).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text update ---------------------------------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value     = "In Translation"
$dede.Range("C2").Value     = "In Translation"

# --- Column width update ---------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F
$zhcn.Columns.Item(3).ColumnWidth     = 12.5   # column C
$dede.Columns.Item(3).ColumnWidth     = 12.5   # column C
